$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing PriceChange / UpDown values for row 5 (previously blank)
$ws.Range("X5").Value = -1.4100040000000149
$ws.Range("Y5").Value = "Down"

# Append a new data row (row 6) with freshly scanned data
$ws.Range("A6").Value = 42647.885428240741
$ws.Range("B6").Value = 8
$ws.Range("C6").Value = "Buy"
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = "Random"
$ws.Range("Q6").Value = 60.94594728999143
$ws.Range("R6").Value = 0
$ws.Range("S6").Value = 0.1189
$ws.Range("T6").Value = 0.008
$ws.Range("U6").Value = 5.99
$ws.Range("V6").Value = "N/A"
$ws.Range("W6").Value = 0

# Match formatting used by the rest of the table (copy number formats from row 5)
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)

$ws.Range("S5").Copy()
$ws.Range("S6").PasteSpecial(-4122)

$ws.Range("T5").Copy()
$ws.Range("T6").PasteSpecial(-4122)

$excel.CutCopyMode = 0
